$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete paragraphs from "Jeong Hoon Choi ..." through the trailing blank
#    paragraph at the very end of the body (original paragraphs 7-18).
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(7)
$lastPara  = $d.Paragraphs.Item($d.Paragraphs.Count)
$rngTail = $d.Range($startPara.Range.Start, $lastPara.Range.End)
$rngTail.Delete()

# ---------------------------------------------------------------------------
# 2) Paragraph 6: "1. " + "Team Formation" -> single run "1.1 "
#    (keep bold formatting, drop the Times New Roman run font override).
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$xmlTeam = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:jc w:val="both"/><w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">1.1 </w:t></w:r>' +
  '</w:p>'
$p6.Range.InsertXML($xmlTeam)

# ---------------------------------------------------------------------------
# 3) Delete paragraph 5 (empty bold paragraph that sat between the hyperlink
#    and the "1. Team Formation" heading).
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Delete()

# ---------------------------------------------------------------------------
# 4) Delete paragraph 4 (the hyperlink-only paragraph).
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Delete()

# ---------------------------------------------------------------------------
# 5) Paragraph 3: drop the inline picture and turn it into the new heading
#    "1 Matrix Multiplication on the CPU" (bold, no explicit run font).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$xmlHeading = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1 Matrix Multiplication on the CPU</w:t></w:r>' +
  '</w:p>'
$p3.Range.InsertXML($xmlHeading)

# ---------------------------------------------------------------------------
# 6) Title paragraph: "Laboratory Assignment 2" -> "...3" (keep the two
#    original runs, only the trailing digit run's text changes).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$xmlTitle = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:jc w:val="center"/><w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr>' +
  '<w:t xml:space="preserve">Laboratory Assignment </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr>' +
  '<w:t>3</w:t></w:r>' +
  '</w:p>'
$p1.Range.InsertXML($xmlTitle)

# ---------------------------------------------------------------------------
# 7) The removed image/hyperlink relationships free up rId7/rId8, so Word
#    renumbers the (now sole) header relationship down to rId7.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
Write-Output ("Header exists: " + $hdr.Exists)

Write-Output "Done"
